# IKD update: GaN CMOS 2026-02-13T23:35Z
# Append 6 new literature records (rows 191-196) to the "Master" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        B="Design and Implementation of a Three-Phase Buck-Boost Split-Source Inverter (BSSI)"; C=2026;
        D="MDPI AG"; E="Electronics"; F="Abdulhussein, Yasameen Sh.; Gün, Ayhan";
        H="10.3390/electronics15040808"; I="https://doi.org/10.3390/electronics15040808"; J="Journal";
        K="Inverter"; L="Experiment"; M="Contacts";
        Q="Design and Implementation of a Three-Phase Buck-Boost Split-Source Inverter (BSSI)"; R="High";
        S="2026-02-13"
    },
    @{
        B="Research on structural reinforcement of AlGaN/GaN HEMT devices under RF stress"; C=2026;
        D="IOP Publishing"; E="Nanotechnology"; F="Liu, Xingjun; Liu, Hongxia; Su, Mengwei; Xing, Dong; Liu, Chang";
        H="10.1088/1361-6528/ae45b0"; I="https://doi.org/10.1088/1361-6528/ae45b0"; J="Journal";
        K="n-FET"; L="TCAD"; M="Contacts";
        Q="Research on structural reinforcement of AlGaN/GaN HEMT devices under RF stress"; R="High";
        S="2026-02-13"
    },
    @{
        B="Design and Implementation of a Three-Phase Buck-Boost Split-Source Inverter (BSSI)"; C=2026;
        D="MDPI AG"; E="Electronics"; F="Abdulhussein, Yasameen Sh.; Gün, Ayhan";
        H="10.3390/electronics15040808"; I="https://doi.org/10.3390/electronics15040808"; J="Journal";
        K="Inverter"; L="Experiment"; M="Contacts";
        Q="Design and Implementation of a Three-Phase Buck-Boost Split-Source Inverter (BSSI)"; R="High";
        S="2026-02-13"
    },
    @{
        B="Robust coordinated fault-tolerant control for aerospace multi-motor synchronous drive systems against inverter fault"; C=2026;
        D="SAGE Publications"; E="Measurement and Control"; F="Han, Xiaodong; Zhang, Dengfeng; Zhao, Li; Lu, Baochun";
        H="10.1177/00202940261419018"; I="https://doi.org/10.1177/00202940261419018"; J="Journal";
        K="Inverter"; L="Experiment"; M="Transport";
        Q="Robust coordinated fault-tolerant control for aerospace multi-motor synchronous drive systems against inverter fault"; R="High";
        S="2026-02-13"
    },
    @{
        B="Design and Implementation of a Three-Phase Buck-Boost Split-Source Inverter (BSSI)"; C=2026;
        D="MDPI AG"; E="Electronics"; F="Abdulhussein, Yasameen Sh.; Gün, Ayhan";
        H="10.3390/electronics15040808"; I="https://doi.org/10.3390/electronics15040808"; J="Journal";
        K="Inverter"; L="Experiment"; M="Contacts";
        Q="Design and Implementation of a Three-Phase Buck-Boost Split-Source Inverter (BSSI)"; R="High";
        S="2026-02-13"
    },
    @{
        B="Skyrmion manipulation and logic gate functionality in transition metal multilayers"; C=2026;
        D="IOP Publishing"; E="Journal of Physics D: Applied Physics"; F="Mukherjee, Tamali; Satya Narayana Murthy, V; Sadhukhan, Banasree";
        H="10.1088/1361-6463/ae45b9"; I="https://doi.org/10.1088/1361-6463/ae45b9"; J="Journal";
        K="Inverter"; L="Experiment"; M="Contacts";
        Q="Skyrmion manipulation and logic gate functionality in transition metal multilayers"; R="High";
        S="2026-02-13"
    }
)

# The sheet consistently leaves columns A, G, N, O, P, T blank (present-but-empty)
# for every existing record. Copying the last populated row down preserves that
# same blank-cell shape for the new rows; we then overwrite the populated fields.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
$startRow = $lastRow + 1
$srcRow = $ws.Range("A$lastRow`:T$lastRow")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $dstRow = $ws.Range("A$r`:T$r")
    $srcRow.Copy($dstRow)

    $rowData = $newRows[$i]
    foreach ($col in $rowData.Keys) {
        $cell = $ws.Range("$col$r")
        $val = $rowData[$col]
        if ($col -eq "S") {
            # Leading apostrophe forces text, so the ISO date string isn't
            # auto-converted to a date serial (matches the source AddedDate
            # column, which stores plain "yyyy-mm-dd" text).
            $cell.Value = "'" + $val
        } else {
            $cell.Value = $val
        }
    }
}
